# Add 2022-Q3 data
# 1) Update the "总计" (summary) sheet: insert a new row right under the
#    header for the 2022-Q3 totals, pushing every existing quarter down.
# 2) Add a brand-new "2022-Q3" worksheet (cloned from "2022-Q2" so it keeps
#    the same column layout/formatting), positioned right after "总计",
#    and fill it with the fund holdings for that quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) 总计 sheet: insert the new 2022-Q3 summary row at row 2
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows(2).Insert()
# Insert() drags a copy of the surrounding formatting along for the ride,
# but not quite the same shape the rest of the table uses - put row 2 back
# in line with its neighbours: column A keeps the bordered/bold "index"
# style, B:D stay plain.
$summary.Range("B2:D2").ClearFormats()
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 18
$summary.Cells.Item(2, 4).Value = 4.59

# Re-number the index column (A) for all the rows that shifted down one slot
for ($r = 3; $r -le 9; $r++) {
    $summary.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------------
# 2) Create the "2022-Q3" worksheet by cloning "2022-Q2" (same columns),
#    inserted right before it, then overwrite with the 2022-Q3 figures.
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q2")
$template.Copy($template, $null)

$newSheet = $wb.Worksheets.Item("2022-Q2 (2)")
$newSheet.Name = "2022-Q3"

# Source data keeps the fund code / amount / weight columns as plain text
# (exactly as they were scraped), so force text formatting before writing
# the values - otherwise "040015" becomes 40015 and "0.00" becomes 0.
$newSheet.Columns("B").NumberFormat = "@"
$newSheet.Columns("D:G").NumberFormat = "@"

$rows = @(
    @("040015", "华安动态灵活配置混合A",       "22.58", "77.56", "4.74", "1.0703", 4),
    @("040001", "华安创新混合",                 "15.29", "72.21", "4.89", "0.7477", 3),
    @("014007", "华安制造升级一年持有混合A",     "12.31", "92.39", "5.04", "0.6204", 3),
    @("010792", "华安成长先锋混合A",             "11.21", "91.59", "5.09", "0.5706", 2),
    @("006154", "华安制造先锋混合A",             "10.65", "93.29", "4.99", "0.5314", 4),
    @("014389", "华安产业动力6个月持有混合A",     "6.79",  "93.58", "5.04", "0.3422", 3),
    @("013619", "华安动态灵活配置混合C",         "4.47",  "77.56", "4.74", "0.2119", 4),
    @("010793", "华安成长先锋混合C",             "3.40",  "91.59", "5.09", "0.1731", 2),
    @("160425", "华安创业板两年定期开放混合",     "1.80",  "93.72", "8.04", "0.1447", 4),
    @("014008", "华安制造升级一年持有混合C",     "0.61",  "92.39", "5.04", "0.0307", 3),
    @("014390", "华安产业动力6个月持有混合C",     "0.59",  "93.58", "5.04", "0.0297", 3),
    @("001675", "江信同福灵活配置混合A",         "0.60",  "91.25", "4.83", "0.0290", 9),
    @("013507", "华安制造先锋混合C",             "0.52",  "93.29", "4.99", "0.0259", 4),
    @("005083", "诺德量化蓝筹增强混合C",         "0.57",  "92.85", "3.36", "0.0192", 5),
    @("004927", "中航军民融合精选混合C",         "0.41",  "72.65", "4.35", "0.0178", 6),
    @("001676", "江信同福灵活配置混合C",         "0.36",  "91.25", "4.83", "0.0174", 9),
    @("004926", "中航军民融合精选混合A",         "0.09",  "72.65", "4.35", "0.0039", 6),
    @("005082", "诺德量化蓝筹增强混合A",         "0.00",  "92.85", "3.36", "__NUM0__", 5)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).Value = $r - 2
    $newSheet.Cells.Item($r, 2).Value = $row[0]
    $newSheet.Cells.Item($r, 3).Value = $row[1]
    $newSheet.Cells.Item($r, 4).Value = $row[2]
    $newSheet.Cells.Item($r, 5).Value = $row[3]
    $newSheet.Cells.Item($r, 6).Value = $row[4]

    if ($row[5] -eq "__NUM0__") {
        # Last row's "持有市值(亿元)" is stored as a real 0, not "0.00" text
        $newSheet.Range("G" + $r).NumberFormat = "General"
        $newSheet.Cells.Item($r, 7).Value = 0
    } else {
        $newSheet.Cells.Item($r, 7).Value = $row[5]
    }

    $newSheet.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# The template had more rows (47) than the new quarter's data (19) - drop the leftovers.
$lastTemplateRow = 47
$firstBlankRow = $rows.Count + 2
$newSheet.Range("A" + $firstBlankRow + ":H" + $lastTemplateRow).EntireRow.Delete()
